# Add one module (fall 2025 semester, 16 weeks) to the course schedule and
# move the view/selection to the newly-added data, mirroring the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Style helpers.
#
# Style indices used by the new rows, matching the cellXfs already present in
# the workbook (style id -> numFmtId/font combination):
#   0 -> General, default font              (no explicit style needed)
#   1 -> "yyyy-mm-dd;@", default font        (e.g. C2, most date cells)
#   3 -> "yyyy-mm-dd", font3                 (e.g. E51, some exam dates)
#   4 -> General, font4                      (e.g. D75, exam labels)
#   5 -> "yyyy-mm-dd", font2                 (e.g. E79, some exam dates)
#   6 -> General, font2                      (brand new style, created below)
#
# The most reliable way to reproduce a cellXf exactly (so Excel reuses the
# existing style entry instead of fabricating a new, slightly different one)
# is to copy an existing cell that already carries that style and then
# overwrite its value.
# ---------------------------------------------------------------------------

$styleSource1 = $ws.Cells.Item(2, 3)    # C2  -> style 1 (yyyy-mm-dd;@)
$styleSource3 = $ws.Cells.Item(51, 5)   # E51 -> style 3 (yyyy-mm-dd, font3)
$styleSource4 = $ws.Cells.Item(75, 4)   # D75 -> style 4 (General, font4)
$styleSource5 = $ws.Cells.Item(79, 5)   # E79 -> style 5 (yyyy-mm-dd, font2)

function Set-CellStyle($cell, $style) {
    if ($style -eq 1) {
        $styleSource1.Copy($cell)
    } elseif ($style -eq 3) {
        $styleSource3.Copy($cell)
    } elseif ($style -eq 4) {
        $styleSource4.Copy($cell)
    } elseif ($style -eq 5) {
        $styleSource5.Copy($cell)
    } elseif ($style -eq 6) {
        # Style 6 (General/font2) does not exist yet in this workbook; giving
        # a plain cell an explicit black font color is what originally
        # produced that exact cellXf combination, so reproduce it the same
        # way here (first use creates it, later uses just reuse it).
        $cell.Font.Color = 0
    }
}

$rows = @(
  @{ row=83; cells=@(@{col=1; style=0; kind="str"; val="fall 2025"}, @{col=2; style=0; kind="num"; val=1}, @{col=3; style=1; kind="num"; val=45894}, @{col=8; style=0; kind="str"; val="Module 1"}, @{col=9; style=0; kind="str"; val="Python Basics (constants, variables, comments, strings, print)"}) },
  @{ row=84; cells=@(@{col=1; style=0; kind="str"; val="fall 2025"}, @{col=2; style=0; kind="num"; val=2}, @{col=3; style=1; kind="num"; val=45901}, @{col=5; style=5; kind="num"; val=45901}, @{col=6; style=6; kind="empty"; val=$null}, @{col=7; style=6; kind="str"; val="Labor Day"}, @{col=8; style=0; kind="str"; val="Module 2"}, @{col=9; style=0; kind="str"; val="Operators and Expressions, functions"}) },
  @{ row=85; cells=@(@{col=1; style=0; kind="str"; val="fall 2025"}, @{col=2; style=0; kind="num"; val=3}, @{col=3; style=1; kind="num"; val=45908}, @{col=8; style=0; kind="str"; val="Module 3"}, @{col=9; style=0; kind="str"; val="Functions, decomposition"}) },
  @{ row=86; cells=@(@{col=1; style=0; kind="str"; val="fall 2025"}, @{col=2; style=0; kind="num"; val=4}, @{col=3; style=1; kind="num"; val=45915}, @{col=8; style=0; kind="str"; val="Module 4"}, @{col=9; style=0; kind="str"; val="Functions, input from user, decomposition"}) },
  @{ row=87; cells=@(@{col=1; style=0; kind="str"; val="fall 2025"}, @{col=2; style=0; kind="num"; val=5}, @{col=3; style=1; kind="num"; val=45922}, @{col=8; style=0; kind="str"; val="Module 5"}, @{col=9; style=0; kind="str"; val="Control Flow (if statements)"}) },
  @{ row=88; cells=@(@{col=1; style=0; kind="str"; val="fall 2025"}, @{col=2; style=0; kind="num"; val=6}, @{col=3; style=1; kind="num"; val=45929}, @{col=8; style=0; kind="str"; val="Module 6"}, @{col=9; style=0; kind="str"; val="Control Flow (while)"}) },
  @{ row=89; cells=@(@{col=1; style=0; kind="str"; val="fall 2025"}, @{col=2; style=0; kind="num"; val=7}, @{col=3; style=1; kind="num"; val=45936}, @{col=4; style=0; kind="str"; val="Midterm 1"}, @{col=5; style=1; kind="num"; val=45937}, @{col=8; style=0; kind="str"; val="Module 7"}, @{col=9; style=0; kind="str"; val="Data Structures (lists)"}) },
  @{ row=90; cells=@(@{col=1; style=0; kind="str"; val="fall 2025"}, @{col=2; style=0; kind="num"; val=8}, @{col=3; style=1; kind="num"; val=45943}, @{col=8; style=0; kind="str"; val="Module 8"}, @{col=9; style=0; kind="str"; val="Control Flow (for loops), mutability, random"}) },
  @{ row=91; cells=@(@{col=1; style=0; kind="str"; val="fall 2025"}, @{col=2; style=0; kind="num"; val=9}, @{col=3; style=1; kind="num"; val=45950}, @{col=8; style=0; kind="str"; val="Module 9"}, @{col=9; style=0; kind="str"; val="Control Flow (for loops), Dictionaries"}) },
  @{ row=92; cells=@(@{col=1; style=0; kind="str"; val="fall 2025"}, @{col=2; style=0; kind="num"; val=10}, @{col=3; style=1; kind="num"; val=45957}, @{col=8; style=0; kind="str"; val="Module 10"}, @{col=9; style=0; kind="str"; val="Files and strings"}) },
  @{ row=93; cells=@(@{col=1; style=0; kind="str"; val="fall 2025"}, @{col=2; style=0; kind="num"; val=11}, @{col=3; style=1; kind="num"; val=45964}, @{col=4; style=4; kind="str"; val="Midterm 2"}, @{col=5; style=1; kind="num"; val=45965}, @{col=8; style=0; kind="str"; val="Module 11"}, @{col=9; style=0; kind="str"; val="Data Structures (tuples)"}) },
  @{ row=94; cells=@(@{col=1; style=0; kind="str"; val="fall 2025"}, @{col=2; style=0; kind="num"; val=12}, @{col=3; style=1; kind="num"; val=45971}, @{col=5; style=3; kind="num"; val=45972}, @{col=7; style=0; kind="str"; val="Veterans Day"}, @{col=8; style=0; kind="str"; val="Module 12"}, @{col=9; style=0; kind="str"; val="2D lists, nested for loops"}) },
  @{ row=95; cells=@(@{col=1; style=0; kind="str"; val="fall 2025"}, @{col=2; style=0; kind="num"; val=13}, @{col=3; style=1; kind="num"; val=45978}, @{col=8; style=0; kind="str"; val="Module 13"}, @{col=9; style=0; kind="str"; val="Data Structures (sets)"}) },
  @{ row=96; cells=@(@{col=1; style=0; kind="str"; val="fall 2025"}, @{col=2; style=0; kind="num"; val=14}, @{col=3; style=1; kind="num"; val=45985}, @{col=5; style=3; kind="num"; val=45988}, @{col=7; style=0; kind="str"; val="Thanksgiving Recess"}, @{col=8; style=0; kind="str"; val="Module 14"}, @{col=9; style=0; kind="str"; val="Mutability"}) },
  @{ row=97; cells=@(@{col=1; style=0; kind="str"; val="fall 2025"}, @{col=2; style=0; kind="num"; val=15}, @{col=3; style=1; kind="num"; val=45992}, @{col=8; style=0; kind="str"; val="Module 15"}, @{col=9; style=0; kind="str"; val="Control Flow + Data Structures"}) },
  @{ row=98; cells=@(@{col=1; style=0; kind="str"; val="fall 2026"}, @{col=2; style=0; kind="num"; val=16}, @{col=3; style=1; kind="num"; val=45999}, @{col=4; style=0; kind="str"; val="Final Exam"}, @{col=5; style=1; kind="num"; val=46001}, @{col=6; style=0; kind="str"; val="TBD"}, @{col=8; style=0; kind="str"; val="Module 16"}, @{col=9; style=0; kind="str"; val="Review, Final Exam"}) },
)

foreach ($rowInfo in $rows) {
    $r = $rowInfo.row
    foreach ($cellInfo in $rowInfo.cells) {
        $cell = $ws.Cells.Item($r, $cellInfo.col)

        # Apply formatting (number format / font) first so that it is in
        # place before (and is not clobbered by) the value assignment.
        Set-CellStyle $cell $cellInfo.style

        if ($cellInfo.kind -eq "str") {
            $cell.Value = $cellInfo.val
        } elseif ($cellInfo.kind -eq "num") {
            $cell.Value = $cellInfo.val
        }
        # "empty" cells (e.g. F84) only need the style applied, no value.
    }
}

# Scroll/select the newly added data, same as the source workbook ends up
# doing: the active cell lands on E98 (the final exam date of the new
# module).
$ws.Range("E98").Select()
